# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (per-fund holdings detail) right
#   before the "总计" (totals) worksheet.
# - Inserts a new leading row into the "总计" worksheet summarizing the
#   2022-Q1 quarter (2 holdings, 0.34 亿元), pushing the existing
#   2021-Q4 / 2020-Q4 rows down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet immediately before "总计".
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet references track *position*, not identity - now that the
# new sheet occupies the old "总计" slot, re-resolve "总计" by name so we
# don't accidentally keep operating on the new sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row (row 1), columns B:H - bold, centered, thin-bordered like the
# other quarter sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: index column (A) + fund code/name/size/position/weight/value/rank.
$rows = @(
    @("003721", "易方达标普信息科技指数（QDII-LOF）美元", "6.31", "93.58", "2.69", "0.1697", 5),
    @("161128", "易方达标普信息科技指数（QDII-LOF）人民币", "6.31", "93.58", "2.69", "0.1697", 5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = 2 + $r
    $idxCell = $newSheet.Cells.Item($rowNum, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $data = $rows[$r]
    # Fund code and the numeric-looking measurements are stored as literal
    # text in the source data (e.g. "003721" keeps its leading zero) - the
    # leading apostrophe forces Excel to keep them as text instead of
    # auto-converting to numbers.
    $newSheet.Cells.Item($rowNum, 2).Value = "'" + $data[0]
    $newSheet.Cells.Item($rowNum, 3).Value = $data[1]
    $newSheet.Cells.Item($rowNum, 4).Value = "'" + $data[2]
    $newSheet.Cells.Item($rowNum, 5).Value = "'" + $data[3]
    $newSheet.Cells.Item($rowNum, 6).Value = "'" + $data[4]
    $newSheet.Cells.Item($rowNum, 7).Value = "'" + $data[5]
    $newSheet.Cells.Item($rowNum, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2) Insert a new leading data row into "总计" for the 2022-Q1 quarter.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
# Excel's row insert copies the formatting of the row above down onto the
# freshly-inserted row; strip that back to the sheet's default (unstyled)
# look before applying the formatting this row actually needs.
$totalSheet.Range("A2:D2").ClearFormats()

$a2 = $totalSheet.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.34

# Renumber the index column (A) for the rows that shifted down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
